# Dicionario de Classes Projeto Ponto.xlsx
# "Employee WorkLoad Edit Option Avaliable for Admins Now"
#
# The 4th sheet (tb_cargas_horarias) documents the "carga horaria" table.
# Two fields already existed (horario_entrada / horario_saida) but their
# "Not Null" restriction was left blank. This change:
#   - fills in the missing "Not Null" restriction + note for those two rows
#   - documents two new fields: horario_comeco_pausa / horario_fim_pausa
#   - leaves the workbook scrolled/focused on that sheet, like the author did

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tb_cargas_horarias")
$ws.Activate()

# --- fill in the previously-empty "Restricoes" / "Observacoes" columns for
#     the two existing rows (horario_entrada / horario_saida) ---
$ws.Range("F10").Value = "Not Null"
$ws.Range("G10").Value = "Este campo não pode ser nulo"

$ws.Range("F11").Value = "Not Null"
$ws.Range("G11").Value = "Este Campo não pode ser nulo"

# --- add two new rows documenting the break-time fields, matching the
#     look/formatting of the existing horario_entrada / horario_saida rows ---
$ws.Range("B10:G10").Copy()
$ws.Range("B12:G12").PasteSpecial(-4122)

$ws.Range("B11:G11").Copy()
$ws.Range("B13:G13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Rows.Item(12).RowHeight = 30.75
$ws.Rows.Item(13).RowHeight = 30.75

$ws.Range("B12").WrapText = $true
$ws.Range("B13").WrapText = $true

$ws.Range("B12").Value = "horario_comeco`n_pausa"
$ws.Range("C12").Value = "Horário que o funcionário`n deve começar sua pausa"
$ws.Range("D12").Value = "time"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""

$ws.Range("B13").Value = "horario_fim`n_pausa"
$ws.Range("C13").Value = "Horário que o funcionário`n deve finalizar sua pausa"
$ws.Range("D13").Value = "time"
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = ""

# --- leave the view the way the author left it: scrolled down so the new
#     rows are visible, with G13 as the active selection ---
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("G13").Select()

# --- bring the tb_cargas_horarias tab into view among the sheet tabs ---
$wb.Windows.Item(1).ScrollWorkbookTabs(3)
